$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51 (hunk 0)
$ws.Range("H51").Value = 3750
$ws.Range("I51").Value = 2916.6667
$ws.Range("J51").Value = 4000
$ws.Range("K51").Value = 2916.6667
$ws.Range("L51").Value = 4000
$ws.Range("M51").Value = -2432.6667
$ws.Range("N51").Value = -4968

# Row 135 (hunk 1)
$ws.Range("H135").Value = 528.44446
$ws.Range("I135").Value = 528.44446
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4756.00014
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2221.00014
$ws.Range("N135").ClearContents()

# Row 137 (hunk 2)
$ws.Range("H137").Value = 1454675.8
$ws.Range("I137").Value = 6214.6665
$ws.Range("J137").Value = 2075444.8
$ws.Range("K137").Value = 18643.9995
$ws.Range("L137").Value = 6226334.4
$ws.Range("M137").Value = -16093.9995
$ws.Range("N137").Value = -6231434.4

# Row 138 (hunk 3)
$ws.Range("H138").Value = 1672.8536
$ws.Range("J138").Value = 3163.8
$ws.Range("L138").Value = 9491.400000000001
$ws.Range("N138").Value = -19771.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 4)
$ws.Range("H32").Value = 3787.97
$ws.Range("I32").Value = 3196.6365
$ws.Range("J32").Value = 8124.4165
$ws.Range("K32").Value = 3196.6365
$ws.Range("L32").Value = 8124.4165
$ws.Range("M32").Value = -2909.6365
$ws.Range("N32").Value = -8698.416499999999

# Row 74 (hunk 5)
$ws.Range("H74").Value = 42438.88
$ws.Range("I74").Value = 78671.62
$ws.Range("J74").Value = 3186.75
$ws.Range("K74").Value = 78671.62
$ws.Range("L74").Value = 3186.75
$ws.Range("M74").Value = -77797.62
$ws.Range("N74").Value = -4934.75

# Row 77 (hunk 6)
$ws.Range("H77").Value = 42438.88
$ws.Range("I77").Value = 78671.62
$ws.Range("J77").Value = 3186.75
$ws.Range("K77").Value = 393358.1
$ws.Range("L77").Value = 15933.75
$ws.Range("M77").Value = -388990.1
$ws.Range("N77").Value = -24669.75

# Row 97 (hunk 7)
$ws.Range("H97").Value = 703.40625
$ws.Range("I97").Value = 688.5714
$ws.Range("J97").Value = 807.25
$ws.Range("K97").Value = 688.5714
$ws.Range("L97").Value = 807.25
$ws.Range("M97").Value = -192.5714
$ws.Range("N97").Value = -1799.25

$ws = $wb.Worksheets.Item("BSM")
# Row 134 (hunk 8)
$ws.Range("H134").Value = 1641.75
$ws.Range("I134").Value = 889.6053000000001
$ws.Range("J134").Value = 4499.9
$ws.Range("K134").Value = 2668.8159
$ws.Range("L134").Value = 13499.7
$ws.Range("M134").Value = -133.8159000000001
$ws.Range("N134").Value = -18569.7

$ws = $wb.Worksheets.Item("CRP")
# Row 25 (hunk 9)
$ws.Range("H25").Value = 8799.6
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

# Row 31 (hunk 10)
$ws.Range("H31").Value = 2379.7715
$ws.Range("I31").Value = 1433.0834
$ws.Range("J31").Value = 4445.273
$ws.Range("K31").Value = 1433.0834
$ws.Range("L31").Value = 4445.273
$ws.Range("M31").Value = -1138.0834
$ws.Range("N31").Value = -5035.273

# Row 34 (hunk 11)
$ws.Range("H34").Value = 2379.7715
$ws.Range("I34").Value = 1433.0834
$ws.Range("J34").Value = 4445.273
$ws.Range("K34").Value = 1433.0834
$ws.Range("L34").Value = 4445.273
$ws.Range("M34").Value = -1231.0834
$ws.Range("N34").Value = -4849.273

# Row 58 (hunk 12)
$ws.Range("H58").Value = 1430538
$ws.Range("I58").Value = 2501589
$ws.Range("J58").Value = 2469.9333
$ws.Range("K58").Value = 2501589
$ws.Range("L58").Value = 2469.9333
$ws.Range("M58").Value = -2501386
$ws.Range("N58").Value = -2875.9333

# Row 105 (hunk 13)
$ws.Range("H105").Value = 2432.1155
$ws.Range("I105").Value = 1695.8334
$ws.Range("K105").Value = 1695.8334
$ws.Range("M105").Value = 51.16660000000002

# Row 136 (hunk 14)
$ws.Range("H136").Value = 1430538
$ws.Range("I136").Value = 2501589
$ws.Range("J136").Value = 2469.9333
$ws.Range("K136").Value = 7504767
$ws.Range("L136").Value = 7409.7999
$ws.Range("M136").Value = -7502217
$ws.Range("N136").Value = -12509.7999

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (hunk 15)
$ws.Range("H5").Value = 1442.2
$ws.Range("I5").Value = 1201.5714
$ws.Range("J5").Value = 1652.75
$ws.Range("K5").Value = 3604.7142
$ws.Range("L5").Value = 4958.25
$ws.Range("M5").Value = -3492.7142
$ws.Range("N5").Value = -5182.25

# Row 32 (hunk 16)
$ws.Range("H32").Value = 25928
$ws.Range("I32").Value = 25049.5
$ws.Range("J32").Value = 26367.25
$ws.Range("K32").Value = 75148.5
$ws.Range("L32").Value = 79101.75
$ws.Range("M32").Value = -74865.5
$ws.Range("N32").Value = -79667.75

# Row 34 (hunk 17)
$ws.Range("H34").Value = 458.6
$ws.Range("I34").Value = 430
$ws.Range("J34").Value = 501.5
$ws.Range("K34").Value = 1290
$ws.Range("L34").Value = 1504.5
$ws.Range("M34").Value = -1206
$ws.Range("N34").Value = -1672.5

# Row 63 (hunk 18)
$ws.Range("H63").Value = 4181.5
$ws.Range("I63").Value = 312
$ws.Range("J63").Value = 5471.3335
$ws.Range("K63").Value = 936
$ws.Range("L63").Value = 16414.0005
$ws.Range("M63").Value = -187
$ws.Range("N63").Value = -17912.0005

# Row 66 (hunk 19)
$ws.Range("H66").Value = 4181.5
$ws.Range("I66").Value = 312
$ws.Range("J66").Value = 5471.3335
$ws.Range("K66").Value = 2808
$ws.Range("L66").Value = 49242.0015
$ws.Range("M66").Value = 936
$ws.Range("N66").Value = -56730.0015

# Row 100 (hunk 20)
$ws.Range("H100").Value = 20667
$ws.Range("J100").Value = 20667
$ws.Range("L100").Value = 62001
$ws.Range("N100").Value = -63623

# Row 135 (hunk 21)
$ws.Range("H135").Value = 1442.2
$ws.Range("I135").Value = 1201.5714
$ws.Range("J135").Value = 1652.75
$ws.Range("K135").Value = 10814.1426
$ws.Range("L135").Value = 14874.75
$ws.Range("M135").Value = -8279.142600000001
$ws.Range("N135").Value = -19944.75

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (hunk 22)
$ws.Range("H102").Value = 1118.3334
$ws.Range("I102").Value = 1134.4231
$ws.Range("J102").Value = 700
$ws.Range("K102").Value = 1134.4231
$ws.Range("L102").Value = 700
$ws.Range("M102").Value = 487.5769
$ws.Range("N102").Value = -3944

# Row 126 (hunk 23)
$ws.Range("H126").Value = 3195.125
$ws.Range("I126").Value = 2224.739
$ws.Range("J126").Value = 5675
$ws.Range("K126").Value = 6674.217000000001
$ws.Range("L126").Value = 17025
$ws.Range("M126").Value = -4204.217000000001
$ws.Range("N126").Value = -21965

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (hunk 24)
$ws.Range("H7").Value = 3985.077
$ws.Range("I7").Value = 4035.9167
$ws.Range("K7").Value = 4035.9167
$ws.Range("M7").Value = -3923.9167

# Row 68 (hunk 25)
$ws.Range("H68").Value = 4243.222
$ws.Range("I68").Value = 4400
$ws.Range("J68").Value = 3929.6667
$ws.Range("K68").Value = 4400
$ws.Range("L68").Value = 3929.6667
$ws.Range("M68").Value = -3651
$ws.Range("N68").Value = -5427.6667

# Row 71 (hunk 26)
$ws.Range("H71").Value = 4243.222
$ws.Range("I71").Value = 4400
$ws.Range("J71").Value = 3929.6667
$ws.Range("K71").Value = 22000
$ws.Range("L71").Value = 19648.3335
$ws.Range("M71").Value = -18256
$ws.Range("N71").Value = -27136.3335

# Row 126 (hunk 27)
$ws.Range("H126").Value = 3985.077
$ws.Range("I126").Value = 4035.9167
$ws.Range("K126").Value = 12107.7501
$ws.Range("M126").Value = -9637.750100000001
